# Car ("汽車") sheet: add capacity column + full metadata columns, and turn
# row 1 into a proper header row (it used to be a stray duplicate of row 2).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: headers -------------------------------------------------------
# B1/D1/E1/F1 already carry header-row style (s=1); just fix their text.
$ws.Cells.Item(1,2).Value = "name"
$ws.Cells.Item(1,4).Value = "owner"
$ws.Cells.Item(1,5).Value = "register_date"
$ws.Cells.Item(1,6).Value = "register_reason"

# C1/G1 exist but are plain numbers today; G1 keeps header style already.
$ws.Cells.Item(1,7).Value = "acquire_value"

# New header cells: copy formatting from an existing header cell, then set text.
$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,3))
$ws.Cells.Item(1,3).Value = "capacity"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,8))
$ws.Cells.Item(1,8).Value = "property_category"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,9))
$ws.Cells.Item(1,9).Value = "category"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,10))
$ws.Cells.Item(1,10).Value = "date"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,11))
$ws.Cells.Item(1,11).Value = "legislator_name"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,12))
$ws.Cells.Item(1,12).Value = "legislator_id"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,13))
$ws.Cells.Item(1,13).Value = "source_file"

$ws.Cells.Item(1,2).Copy($ws.Cells.Item(1,14))
$ws.Cells.Item(1,14).Value = "index"

# --- Rows 2-3: metadata columns H:N, copying format from an existing data cell ---
$rows = @(2, 3)
foreach ($r in $rows) {
    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,8))
    $ws.Cells.Item($r,8).Value = "land"

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,9))
    $ws.Cells.Item($r,9).Value = "normal"

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,10))
    $ws.Cells.Item($r,10).Value = "2011-12-20"

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,11))
    $ws.Cells.Item($r,11).Value = "黃偉哲"

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,12))
    $ws.Cells.Item($r,12).Value = 1367

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,13))
    $ws.Cells.Item($r,13).Value = "tmp85f1"

    $ws.Cells.Item(2,2).Copy($ws.Cells.Item($r,14))
    $ws.Cells.Item($r,14).Value = $r + 28
}
